$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 62, shifting existing rows 62-102 down to 63-103
$ws.Rows.Item(62).Insert()

# Populate the new row 62 with the new record's data
$ws.Cells.Item(62, 1).Value = 5
$ws.Cells.Item(62, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(62, 3).Value = "Maule"
$ws.Cells.Item(62, 4).Value = 44587
$ws.Cells.Item(62, 5).Value = 7
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100108
$ws.Cells.Item(62, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(62, 9).Value = 100108002
$ws.Cells.Item(62, 10).Value = "Mango"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 300
$ws.Cells.Item(62, 14).Value = 7000
$ws.Cells.Item(62, 15).Value = 7000
$ws.Cells.Item(62, 16).Value = 7000
$ws.Cells.Item(62, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(62, 18).Value = "Perú"
$ws.Cells.Item(62, 19).Value = 1750
$ws.Cells.Item(62, 20).Value = 4
